# Refresh the cryptocurrency prices / 1h-volume table (and the swapped
# EnergySwap / PancakeSwap row order) with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.244.52'
$ws.Range('E2').Value = '  +13.64%  '

$ws.Range('D3').Value = '1.675.73'
$ws.Range('E3').Value = '  +8.16%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.34%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '309.04'
$ws.Range('E5').Value = '  +9.18%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9981'
$ws.Range('E6').Value = '  +3.02%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3739'
$ws.Range('E7').Value = '  +3.31%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3440'
$ws.Range('E8').Value = '  +7.17%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.78'
$ws.Range('E9').Value = '  +16.62%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.181'
$ws.Range('E10').Value = '  +6.27%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07310'
$ws.Range('E11').Value = '  +5.34%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.15%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.49'
$ws.Range('E13').Value = '  +8.12%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.110'
$ws.Range('E14').Value = '  +6.59%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.774'
$ws.Range('E15').Value = '  +5.28%  '

$ws.Range('D16').Value = '1.674.88'
$ws.Range('E16').Value = '  +8.42%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001111'
$ws.Range('E17').Value = '  +5.41%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9985'
$ws.Range('E18').Value = '  +3.05%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06723'
$ws.Range('E19').Value = '  +9.29%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '81.99'
$ws.Range('E20').Value = '  +12.01%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.47'
$ws.Range('E21').Value = '  +7.61%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.150'
$ws.Range('E22').Value = '  +6.86%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.01'
$ws.Range('E23').Value = '  +5.13%  '

$ws.Range('D24').Value = '24.216.23'
$ws.Range('E24').Value = '  +13.40%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.413'
$ws.Range('E25').Value = '  +3.92%  '

$ws.Range('E26').Value = '  -9.33%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.669'
$ws.Range('E27').Value = '  +16.91%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '151.75'
$ws.Range('E28').Value = '  +2.62%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.57'
$ws.Range('E29').Value = '  +9.28%  '

$ws.Range('D30').Value = '1.857.35'
$ws.Range('E30').Value = '  +8.27%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '127.24'
$ws.Range('E31').Value = '  +7.16%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.453'
$ws.Range('E32').Value = '  +22.29%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.161'
$ws.Range('E33').Value = '  +2.77%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9992'
$ws.Range('E34').Value = '  +14.31%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.772'
$ws.Range('E35').Value = '  +15.46%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08511'
$ws.Range('E36').Value = '  +5.57%  '

$ws.Range('E37').Value = '  +16.74%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06462'
$ws.Range('E38').Value = '  +10.08%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.381'
$ws.Range('E39').Value = '  +7.61%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.863'
$ws.Range('E40').Value = '  +11.26%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.02356'
$ws.Range('E41').Value = '  +10.75%  '

$ws.Range('E42').Value = '  +5.50%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.2146'
$ws.Range('E43').Value = '  +7.19%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6187'
$ws.Range('E44').Value = '  +12.01%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9973'
$ws.Range('E45').Value = '  +3.04%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.805'
$ws.Range('E46').Value = '  +6.22%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '13.24'
$ws.Range('E47').Value = '  +5.61%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5962'
$ws.Range('E48').Value = '  +8.02%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '127.56'
$ws.Range('E49').Value = '  +4.31%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.035'
$ws.Range('E50').Value = '  +7.87%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07172'
$ws.Range('E51').Value = '  +8.23%  '
